$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains exact text representation (e.g. trailing
# zeros, thousands-separated formatting) instead of being auto-converted to
# numbers by Excel when values look numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.055.83"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "3.794.98"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "601.07"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "165.02"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("D10").Value = "0.451"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("E11").Value = "  +2.80%  "

$ws.Range("D12").Value = "0.0000249"
$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("D13").Value = "35.80"
$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").Value = "4.433.26"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "3.787.36"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "68.106.51"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "18.41"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").Value = "461.47"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").Value = "9.71"
$ws.Range("E21").Value = "  -2.35%  "

$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").Value = "0.0000148"
$ws.Range("E23").Value = "  -4.45%  "

$ws.Range("D24").Value = "83.05"
$ws.Range("E24").Value = "  -0.66%  "

$ws.Range("D25").Value = "12.04"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("D26").Value = "2.11"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "3.945.63"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("E30").Value = "  -5.19%  "

$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("D32").Value = "7.34"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").Value = "29.35"
$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D35").Value = "9.04"
$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("D36").Value = "0.0997"
$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.139"
$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  -3.25%  "

$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("D40").Value = "0.986"
$ws.Range("E40").Value = "  -1.55%  "

$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "47.47"
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").Value = "43.34"
$ws.Range("E45").Value = "  -1.04%  "

$ws.Range("D46").Value = "151.52"
$ws.Range("E46").Value = "  +1.71%  "

$ws.Range("D47").Value = "8.37"
$ws.Range("E47").Value = "  +0.22%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "397.67"
$ws.Range("E48").Value = "  -0.67%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.87"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").Value = "1.35"
$ws.Range("E50").Value = "  +2.61%  "

$ws.Range("D51").Value = "26.50"
$ws.Range("E51").Value = "  -0.94%  "
